$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "91.231.14"
$ws.Cells.Item(2,5).Value = "  +3.91%  "

$ws.Cells.Item(3,4).Value = "3.100.85"
$ws.Cells.Item(3,5).Value = "  -0.18%  "

$c = $ws.Cells.Item(4,4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Cells.Item(4,5).Value = "  -0.14%  "

$c = $ws.Cells.Item(5,4)
$c.NumberFormat = "@"
$c.Value = "218.69"
$c.Style = "Normal"
$ws.Cells.Item(5,5).Value = "  +2.96%  "

$c = $ws.Cells.Item(6,4)
$c.NumberFormat = "@"
$c.Value = "618.81"
$c.Style = "Normal"
$ws.Cells.Item(6,5).Value = "  -2.39%  "

$c = $ws.Cells.Item(7,4)
$c.NumberFormat = "@"
$c.Value = "0.379"
$c.Style = "Normal"
$ws.Cells.Item(7,5).Value = "  -0.35%  "

$c = $ws.Cells.Item(8,4)
$c.NumberFormat = "@"
$c.Value = "0.908"
$c.Style = "Normal"
$ws.Cells.Item(8,5).Value = "  +8.19%  "

$ws.Cells.Item(10,4).Value = "3.097.92"
$ws.Cells.Item(10,5).Value = "  -0.19%  "

$c = $ws.Cells.Item(11,4)
$c.NumberFormat = "@"
$c.Value = "0.674"
$c.Style = "Normal"
$ws.Cells.Item(11,5).Value = "  +13.08%  "

$c = $ws.Cells.Item(12,4)
$c.NumberFormat = "@"
$c.Value = "0.189"
$c.Style = "Normal"
$ws.Cells.Item(12,5).Value = "  +5.94%  "

$c = $ws.Cells.Item(13,4)
$c.NumberFormat = "@"
$c.Value = "0.0000256"
$c.Style = "Normal"
$ws.Cells.Item(13,5).Value = "  +4.71%  "

$ws.Cells.Item(14,4).Value = "91.015.38"
$ws.Cells.Item(14,5).Value = "  +3.74%  "

$ws.Cells.Item(15,5).Value = "  -0.02%  "

$c = $ws.Cells.Item(16,4)
$c.NumberFormat = "@"
$c.Value = "33.07"
$c.Style = "Normal"
$ws.Cells.Item(16,5).Value = "  +3.80%  "

$ws.Cells.Item(17,4).Value = "3.656.25"
$ws.Cells.Item(17,5).Value = "  -0.78%  "

$ws.Cells.Item(18,4).Value = "3.097.70"
$ws.Cells.Item(18,5).Value = "  -0.16%  "

$c = $ws.Cells.Item(19,4)
$c.NumberFormat = "@"
$c.Value = "3.57"
$c.Style = "Normal"
$ws.Cells.Item(19,5).Value = "  +6.27%  "

$c = $ws.Cells.Item(20,4)
$c.NumberFormat = "@"
$c.Value = "0.0000224"
$c.Style = "Normal"
$ws.Cells.Item(20,5).Value = "  +5.35%  "

$c = $ws.Cells.Item(21,4)
$c.NumberFormat = "@"
$c.Value = "13.80"
$c.Style = "Normal"
$ws.Cells.Item(21,5).Value = "  +4.08%  "

$c = $ws.Cells.Item(22,4)
$c.NumberFormat = "@"
$c.Value = "434.44"
$c.Style = "Normal"
$ws.Cells.Item(22,5).Value = "  +2.70%  "

$c = $ws.Cells.Item(23,4)
$c.NumberFormat = "@"
$c.Value = "8.50"
$c.Style = "Normal"
$ws.Cells.Item(23,5).Value = "  +1.00%  "

$c = $ws.Cells.Item(24,4)
$c.NumberFormat = "@"
$c.Value = "5.12"
$c.Style = "Normal"
$ws.Cells.Item(24,5).Value = "  +4.40%  "

$ws.Cells.Item(25,5).Value = "  +2.53%  "

$c = $ws.Cells.Item(26,4)
$c.NumberFormat = "@"
$c.Value = "84.08"
$c.Style = "Normal"
$ws.Cells.Item(26,5).Value = "  +1.07%  "

$c = $ws.Cells.Item(27,4)
$c.NumberFormat = "@"
$c.Value = "11.86"
$c.Style = "Normal"
$ws.Cells.Item(27,5).Value = "  +4.35%  "

$ws.Cells.Item(28,4).Value = "3.248.84"
$ws.Cells.Item(28,5).Value = "  -0.79%  "

$ws.Cells.Item(29,5).Value = "  -0.16%  "

$ws.Cells.Item(30,5).Value = "  +8.14%  "

$ws.Cells.Item(31,5).Value = "  +1.19%  "

$c = $ws.Cells.Item(32,4)
$c.NumberFormat = "@"
$c.Value = "8.71"
$c.Style = "Normal"
$ws.Cells.Item(32,5).Value = "  +7.23%  "

$c = $ws.Cells.Item(33,4)
$c.NumberFormat = "@"
$c.Value = "3.93"
$c.Style = "Normal"
$ws.Cells.Item(33,5).Value = "  +2.70%  "

$c = $ws.Cells.Item(34,4)
$c.NumberFormat = "@"
$c.Value = "520.02"
$c.Style = "Normal"
$ws.Cells.Item(34,5).Value = "  +3.49%  "

$c = $ws.Cells.Item(35,4)
$c.NumberFormat = "@"
$c.Value = "7.04"
$c.Style = "Normal"
$ws.Cells.Item(35,5).Value = "  +4.18%  "

$ws.Cells.Item(38,5).Value = "  +0.91%  "

$ws.Cells.Item(39,5).Value = "  +2.44%  "

$c = $ws.Cells.Item(40,4)
$c.NumberFormat = "@"
$c.Value = "22.30"
$c.Style = "Normal"
$ws.Cells.Item(40,5).Value = "  +0.61%  "

$ws.Cells.Item(41,5).Value = "  -0.14%  "

$ws.Cells.Item(42,5).Value = "  -0.01%  "

$ws.Cells.Item(43,5).Value = "  +1.61%  "

$c = $ws.Cells.Item(44,4)
$c.NumberFormat = "@"
$c.Value = "0.371"
$c.Style = "Normal"
$ws.Cells.Item(44,5).Value = "  +1.40%  "

$ws.Cells.Item(45,5).Value = "  +2.26%  "

$c = $ws.Cells.Item(46,4)
$c.NumberFormat = "@"
$c.Value = "0.0724"
$c.Style = "Normal"
$ws.Cells.Item(46,5).Value = "  +9.66%  "

$ws.Cells.Item(49,5).Value = "  +14.57%  "

$c = $ws.Cells.Item(50,4)
$c.NumberFormat = "@"
$c.Value = "4.20"
$c.Style = "Normal"
$ws.Cells.Item(50,5).Value = "  +6.39%  "


# Row 36: Fetch.AI -> Kaspa
$ws.Cells.Item(36,2).Value = "Kaspa"
$ws.Cells.Item(36,3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Cells.Item(36,4)
$c.NumberFormat = "@"
$c.Value = "0.140"
$c.Style = "Normal"
$ws.Cells.Item(36,5).Value = "  -5.15%  "

# Row 37: Kaspa -> Fetch.AI
$ws.Cells.Item(37,2).Value = "Fetch.AI"
$ws.Cells.Item(37,3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Cells.Item(37,4)
$c.NumberFormat = "@"
$c.Value = "1.29"
$c.Style = "Normal"
$ws.Cells.Item(37,5).Value = "  +2.02%  "

# Row 47: OKB -> Monero
$ws.Cells.Item(47,2).Value = "Monero"
$ws.Cells.Item(47,3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Cells.Item(47,4)
$c.NumberFormat = "@"
$c.Value = "142.73"
$c.Style = "Normal"
$ws.Cells.Item(47,5).Value = "  -2.63%  "

# Row 48: Monero -> OKB
$ws.Cells.Item(48,2).Value = "OKB"
$ws.Cells.Item(48,3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Cells.Item(48,4)
$c.NumberFormat = "@"
$c.Value = "43.57"
$c.Style = "Normal"
$ws.Cells.Item(48,5).Value = "  -0.79%  "

# Row 51: Aave -> ImmutableX
$ws.Cells.Item(51,2).Value = "ImmutableX"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Cells.Item(51,4)
$c.NumberFormat = "@"
$c.Value = "1.25"
$c.Style = "Normal"
$ws.Cells.Item(51,5).Value = "  +5.31%  "
